$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F is "dSF". Update the values per the repulled/recalculated data.
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = -5
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = -12
$ws.Range("F11").Value = -4
$ws.Range("F12").Value = -2
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = 4
